$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.256.10"
$ws.Range("E2").Value = "  -2.89%  "
$ws.Range("D3").Value = "1.549.86"
$ws.Range("E3").Value = "  -4.94%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'206.59"
$ws.Range("E5").Value = "  -3.57%  "
$ws.Range("D6").Value = "'1.01"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.477"
$ws.Range("E7").Value = "  -5.24%  "
$ws.Range("D8").Value = "'0.0607"
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("D10").Value = "'17.63"
$ws.Range("E10").Value = "  -5.22%  "
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").Value = "1.766.20"
$ws.Range("E12").Value = "  -4.92%  "
$ws.Range("D13").Value = "1.552.14"
$ws.Range("E13").Value = "  -4.83%  "
$ws.Range("D14").Value = "'3.97"
$ws.Range("E14").Value = "  -4.99%  "
$ws.Range("D15").Value = "'0.503"
$ws.Range("E15").Value = "  -4.89%  "
$ws.Range("D16").Value = "25.245.80"
$ws.Range("E16").Value = "  -2.95%  "
$ws.Range("D17").Value = "0.0₃0703"
$ws.Range("E17").Value = "  -5.02%  "
$ws.Range("D18").Value = "'58.45"
$ws.Range("E18").Value = "  -4.92%  "
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "'185.13"
$ws.Range("E21").Value = "  -3.95%  "
$ws.Range("D22").Value = "'9.22"
$ws.Range("E22").Value = "  -3.43%  "
$ws.Range("D23").Value = "'5.82"
$ws.Range("E24").Value = "  -4.21%  "
$ws.Range("D25").Value = "'1.01"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "'138.75"
$ws.Range("E26").Value = "  -3.80%  "
$ws.Range("E27").Value = "  -5.50%  "
$ws.Range("D28").Value = "'14.79"
$ws.Range("E28").Value = "  -3.15%  "
$ws.Range("D29").Value = "'6.37"
$ws.Range("E29").Value = "  -5.53%  "
$ws.Range("E30").Value = "  -6.66%  "
$ws.Range("D31").Value = "'0.0462"
$ws.Range("E31").Value = "  -4.54%  "
$ws.Range("D32").Value = "'3.01"
$ws.Range("E32").Value = "  -3.72%  "
$ws.Range("E33").Value = "  -5.51%  "
$ws.Range("D34").Value = "'1.44"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("E35").Value = "  -3.92%  "
$ws.Range("D36").Value = "1.079.49"
$ws.Range("E36").Value = "  -3.76%  "
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("E38").Value = "  -3.20%  "
$ws.Range("E39").Value = "  -5.71%  "
$ws.Range("E40").Value = "  -7.95%  "
$ws.Range("D41").Value = "'0.759"
$ws.Range("E41").Value = "  -10.90%  "
$ws.Range("D42").Value = "'0.798"
$ws.Range("E42").Value = "  +3.85%  "
$ws.Range("D43").Value = "'92.67"
$ws.Range("E43").Value = "  -5.63%  "
$ws.Range("D44").Value = "'5.03"
$ws.Range("E44").Value = "  -2.42%  "
$ws.Range("D45").Value = "1.681.84"
$ws.Range("E45").Value = "  -4.85%  "
$ws.Range("E46").Value = "  -2.40%  "
$ws.Range("D47").Value = "'1.45"
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("D48").Value = "'52.16"
$ws.Range("E48").Value = "  -4.31%  "
$ws.Range("E49").Value = "  -5.20%  "
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("E51").Value = "  -2.10%  "
